# Analysis of predicate frequency.
# Adds a "working" sheet that pulls out the value/z-score column pairs for
# 4 specific "Frequency rank" words from "main record" (ranks 2, 4, 27, 306),
# then computes, per author, the total predicate count across those 4 ranks,
# its z-score, the per-rank relative frequency (rank count / total preds) and
# that frequency's z-score.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, placed right after "main record".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "working"

# ---- Pull the 4 rank blocks (rank 2, 4, 27, 306) out of "main record" ----
# main record columns D:E = rank 2, H:I = rank 4, BB:BC = rank 27, WN:WO = rank 306
$ws1.Range("A1:A15").Copy($ws2.Range("A1"))
$ws1.Range("D1:E15").Copy($ws2.Range("B1"))
$ws1.Range("H1:I15").Copy($ws2.Range("D1"))
$ws1.Range("BB1:BC15").Copy($ws2.Range("F1"))
$ws1.Range("WN1:WO15").Copy($ws2.Range("H1"))

# ---- Second block: totals / relative frequency / z-scores ----
$ws2.Range("E17").Value = "sWord / total preds"

$ws2.Range("A18").Value = "total preds"
$ws2.Range("E18").Value = "#-coord-c-pred_co-v"
$ws2.Range("H18").Value = "#-coord-c-coord-c-pred_co-v"
$ws2.Range("K18").Value = "#-coord-c-coord-c-coord-c-pred_co-v"

$ws2.Range("B19").Value = "frequency"
$ws2.Range("C19").Value = "z-score"
$ws2.Range("E19").Value = "frequency"
$ws2.Range("F19").Value = "z-score"
$ws2.Range("H19").Value = "frequency"
$ws2.Range("I19").Value = "z-score"

# authors, reused from rows 3:15
$ws2.Range("A3:A15").Copy($ws2.Range("A20"))

$ws2.Range("B20").Formula = "=B3+D3+F3+H3"
$ws2.Range("C20").Formula = "=STANDARDIZE(B20,$B$34,4)"
$ws2.Range("E20").Formula = "=D3/B20"
$ws2.Range("F20").Formula = "=STANDARDIZE(E20,$E$34,$E$35)"
$ws2.Range("H20").Formula = "=F3/B20"
$ws2.Range("I20").Formula = "=STANDARDIZE(H20,$H$34,$H$35)"
$ws2.Range("K20").Formula = "=H3/B20"
$ws2.Range("L20").Formula = "=STANDARDIZE(K20,$K$34,$K$35)"
$ws2.Range("B21").Formula = "=B4+D4+F4+H4"
$ws2.Range("C21").Formula = "=STANDARDIZE(B21,$B$34,4)"
$ws2.Range("E21").Formula = "=D4/B21"
$ws2.Range("F21").Formula = "=STANDARDIZE(E21,$E$34,$E$35)"
$ws2.Range("H21").Formula = "=F4/B21"
$ws2.Range("I21").Formula = "=STANDARDIZE(H21,$H$34,$H$35)"
$ws2.Range("K21").Formula = "=H4/B21"
$ws2.Range("L21").Formula = "=STANDARDIZE(K21,$K$34,$K$35)"
$ws2.Range("B22").Formula = "=B5+D5+F5+H5"
$ws2.Range("C22").Formula = "=STANDARDIZE(B22,$B$34,4)"
$ws2.Range("E22").Formula = "=D5/B22"
$ws2.Range("F22").Formula = "=STANDARDIZE(E22,$E$34,$E$35)"
$ws2.Range("H22").Formula = "=F5/B22"
$ws2.Range("I22").Formula = "=STANDARDIZE(H22,$H$34,$H$35)"
$ws2.Range("K22").Formula = "=H5/B22"
$ws2.Range("L22").Formula = "=STANDARDIZE(K22,$K$34,$K$35)"
$ws2.Range("B23").Formula = "=B6+D6+F6+H6"
$ws2.Range("C23").Formula = "=STANDARDIZE(B23,$B$34,4)"
$ws2.Range("E23").Formula = "=D6/B23"
$ws2.Range("F23").Formula = "=STANDARDIZE(E23,$E$34,$E$35)"
$ws2.Range("H23").Formula = "=F6/B23"
$ws2.Range("I23").Formula = "=STANDARDIZE(H23,$H$34,$H$35)"
$ws2.Range("K23").Formula = "=H6/B23"
$ws2.Range("L23").Formula = "=STANDARDIZE(K23,$K$34,$K$35)"
$ws2.Range("B24").Formula = "=B7+D7+F7+H7"
$ws2.Range("C24").Formula = "=STANDARDIZE(B24,$B$34,4)"
$ws2.Range("E24").Formula = "=D7/B24"
$ws2.Range("F24").Formula = "=STANDARDIZE(E24,$E$34,$E$35)"
$ws2.Range("H24").Formula = "=F7/B24"
$ws2.Range("I24").Formula = "=STANDARDIZE(H24,$H$34,$H$35)"
$ws2.Range("K24").Formula = "=H7/B24"
$ws2.Range("L24").Formula = "=STANDARDIZE(K24,$K$34,$K$35)"
$ws2.Range("B25").Formula = "=B8+D8+F8+H8"
$ws2.Range("C25").Formula = "=STANDARDIZE(B25,$B$34,4)"
$ws2.Range("E25").Formula = "=D8/B25"
$ws2.Range("F25").Formula = "=STANDARDIZE(E25,$E$34,$E$35)"
$ws2.Range("H25").Formula = "=F8/B25"
$ws2.Range("I25").Formula = "=STANDARDIZE(H25,$H$34,$H$35)"
$ws2.Range("K25").Formula = "=H8/B25"
$ws2.Range("L25").Formula = "=STANDARDIZE(K25,$K$34,$K$35)"
$ws2.Range("B26").Formula = "=B9+D9+F9+H9"
$ws2.Range("C26").Formula = "=STANDARDIZE(B26,$B$34,4)"
$ws2.Range("E26").Formula = "=D9/B26"
$ws2.Range("F26").Formula = "=STANDARDIZE(E26,$E$34,$E$35)"
$ws2.Range("H26").Formula = "=F9/B26"
$ws2.Range("I26").Formula = "=STANDARDIZE(H26,$H$34,$H$35)"
$ws2.Range("K26").Formula = "=H9/B26"
$ws2.Range("L26").Formula = "=STANDARDIZE(K26,$K$34,$K$35)"
$ws2.Range("B27").Formula = "=B10+D10+F10+H10"
$ws2.Range("C27").Formula = "=STANDARDIZE(B27,$B$34,4)"
$ws2.Range("E27").Formula = "=D10/B27"
$ws2.Range("F27").Formula = "=STANDARDIZE(E27,$E$34,$E$35)"
$ws2.Range("H27").Formula = "=F10/B27"
$ws2.Range("I27").Formula = "=STANDARDIZE(H27,$H$34,$H$35)"
$ws2.Range("K27").Formula = "=H10/B27"
$ws2.Range("L27").Formula = "=STANDARDIZE(K27,$K$34,$K$35)"
$ws2.Range("B28").Formula = "=B11+D11+F11+H11"
$ws2.Range("C28").Formula = "=STANDARDIZE(B28,$B$34,4)"
$ws2.Range("E28").Formula = "=D11/B28"
$ws2.Range("F28").Formula = "=STANDARDIZE(E28,$E$34,$E$35)"
$ws2.Range("H28").Formula = "=F11/B28"
$ws2.Range("I28").Formula = "=STANDARDIZE(H28,$H$34,$H$35)"
$ws2.Range("K28").Formula = "=H11/B28"
$ws2.Range("L28").Formula = "=STANDARDIZE(K28,$K$34,$K$35)"
$ws2.Range("B29").Formula = "=B12+D12+F12+H12"
$ws2.Range("C29").Formula = "=STANDARDIZE(B29,$B$34,4)"
$ws2.Range("E29").Formula = "=D12/B29"
$ws2.Range("F29").Formula = "=STANDARDIZE(E29,$E$34,$E$35)"
$ws2.Range("H29").Formula = "=F12/B29"
$ws2.Range("I29").Formula = "=STANDARDIZE(H29,$H$34,$H$35)"
$ws2.Range("K29").Formula = "=H12/B29"
$ws2.Range("L29").Formula = "=STANDARDIZE(K29,$K$34,$K$35)"
$ws2.Range("B30").Formula = "=B13+D13+F13+H13"
$ws2.Range("C30").Formula = "=STANDARDIZE(B30,$B$34,4)"
$ws2.Range("E30").Formula = "=D13/B30"
$ws2.Range("F30").Formula = "=STANDARDIZE(E30,$E$34,$E$35)"
$ws2.Range("H30").Formula = "=F13/B30"
$ws2.Range("I30").Formula = "=STANDARDIZE(H30,$H$34,$H$35)"
$ws2.Range("K30").Formula = "=H13/B30"
$ws2.Range("L30").Formula = "=STANDARDIZE(K30,$K$34,$K$35)"
$ws2.Range("B31").Formula = "=B14+D14+F14+H14"
$ws2.Range("C31").Formula = "=STANDARDIZE(B31,$B$34,4)"
$ws2.Range("E31").Formula = "=D14/B31"
$ws2.Range("F31").Formula = "=STANDARDIZE(E31,$E$34,$E$35)"
$ws2.Range("H31").Formula = "=F14/B31"
$ws2.Range("I31").Formula = "=STANDARDIZE(H31,$H$34,$H$35)"
$ws2.Range("K31").Formula = "=H14/B31"
$ws2.Range("L31").Formula = "=STANDARDIZE(K31,$K$34,$K$35)"
$ws2.Range("B32").Formula = "=B15+D15+F15+H15"
$ws2.Range("C32").Formula = "=STANDARDIZE(B32,$B$34,4)"
$ws2.Range("E32").Formula = "=D15/B32"
$ws2.Range("F32").Formula = "=STANDARDIZE(E32,$E$34,$E$35)"
$ws2.Range("H32").Formula = "=F15/B32"
$ws2.Range("I32").Formula = "=STANDARDIZE(H32,$H$34,$H$35)"
$ws2.Range("K32").Formula = "=H15/B32"
$ws2.Range("L32").Formula = "=STANDARDIZE(K32,$K$34,$K$35)"

# ---- avg / sd footer ----
$ws2.Range("A34").Value = "avg"
$ws2.Range("B34").Formula = "=AVERAGE(B20:B32)"
$ws2.Range("E34").Formula = "=AVERAGE(E20:E32)"
$ws2.Range("H34").Formula = "=AVERAGE(H20:H32)"
$ws2.Range("K34").Formula = "=AVERAGE(K20:K32)"

$ws2.Range("A35").Value = "sd"
$ws2.Range("B35").Formula = "=STDEV.S(B20:B32)"
$ws2.Range("E35").Formula = "=STDEV.S(E20:E32)"
$ws2.Range("H35").Formula = "=STDEV.S(H20:H32)"
$ws2.Range("K35").Formula = "=STDEV.S(K20:K32)"

# ---- column widths on the new sheet (bestFit, matches "main record") ----
$ws2.Columns.Item(1).ColumnWidth = 19.7265625
$ws2.Columns.Item(2).ColumnWidth = 11.81640625
$ws2.Columns.Item(3).ColumnWidth = 12.453125
$ws2.Columns.Item(4).ColumnWidth = 17.90625
$ws2.Columns.Item(5).ColumnWidth = 17.90625
$ws2.Columns.Item(6).ColumnWidth = 24.90625
$ws2.Columns.Item(7).ColumnWidth = 24.90625
$ws2.Columns.Item(8).ColumnWidth = 32
$ws2.Columns.Item(9).ColumnWidth = 32
$ws2.Columns.Item(11).ColumnWidth = 33.1796875

# ---- view state: freeze panes + selection, both sheets ----
$ws1.Application.ActiveWindow.SplitRow = 2
$ws1.Application.ActiveWindow.SplitColumn = 1
$ws1.Range("WJ3").Select()
$ws1.Range("WN1:WO15").Select()

$ws2.Range("B6").Select()
$ws2.Application.ActiveWindow.FreezePanes = $true
$ws2.Range("G27").Select()
